# Generate Report for Handback
# Updates the localization-status workbook to reflect that the handback
# has been generated/synced: the Overview status moves from "Ready for
# handoff" to "Handed back: in sync with en-US", the handback datetimes
# are refreshed, and the stale "version mismatch" error details are
# cleared now that the handback is current.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# Overview sheet: Status columns for zh-cn (E2) and de-de (F2)
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus

# zh-cn sheet: Status (C2), Latest Handback DateTime (K2), Error Detail (P2)
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("K2").Value = "2016-08-29 02:49:02"
$zhcn.Range("P2").Value = ""

# de-de sheet: Status (C2), Latest Handback DateTime (K2), Error Detail (P2)
$dede.Range("C2").Value = $newStatus
$dede.Range("K2").Value = "2016-08-29 02:49:11"
$dede.Range("P2").Value = ""
